# "Code cleanup in progress"
# The sheet was still carrying an ad-hoc month label ("Jul 2022") left
# over from when the report was first created. Rename it to something
# generic that describes the sheet's actual contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Employees"
